$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows above current row 2 to make room for the two new
# feedback entries (157844 - Mart Minas, and 157821 - Cocamar).
$ws.Range("A2:A5").EntireRow.Insert()

# New data for the inserted rows.
$newRows = @(
    @(157844, "Mart Minas, 19/08/2024", "backlog", "Acompanhamento de clientes", "[]", "Cliente solicitou ajustes ou refação", "2024-08-19", "Mart Minas"),
    @(157844, "Mart Minas, 19/08/2024", "backlog", "Acompanhamento de clientes", "[]", "Entregas feitas conforme planejado", "2024-08-19", "Mart Minas"),
    @(157821, "Cocamar, 15/07/2024", "backlog", "Acompanhamento de clientes", "[]", "Cliente pediu proposta", "2024-07-15", "Cocamar"),
    @(157821, "Cocamar, 15/07/2024", "backlog", "Acompanhamento de clientes", "[]", "Resolveu problema", "2024-07-15", "Cocamar")
)

# Force column G (date-looking text like "2024-08-19") to be stored as
# plain text, matching the original inlineStr cells rather than being
# auto-converted into a date serial number.
$ws.Range("G2:G5").NumberFormat = "@"

$r = 2
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}
